# Updated data to reflect new requirement separation
#
# The "Terms Typically Offered" column (D) needs three new columns inserted
# before it: Corequisites, Concurrent, Recommended. Every data row gets "NA"
# in those three new columns, and the old "Terms Typically Offered" values
# shift right along with the column (D -> G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at D, E, F - this pushes the existing
# "Terms Typically Offered" column (and its data) from D to G.
$ws.Columns("D:F").Insert()

# New header row labels for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Every data row (2-46) gets "NA" for the three new columns.
$ws.Range("D2:F46").Value = "NA"
